# cryptos.xlsx refresh: GitHub Actions symbol-list update
# (Fri Dec 16 14:01:04 UTC 2022)
#
# Updates coin prices (col D) and the scrape hour (col G) for rows 2-51,
# plus the two "Worst/Best in 24h" call-out suffixes that moved from the
# One row (E18) to the SpecialPowerGold row (E51).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Hora columns store numeric-looking values as plain TEXT in this
# sheet. A bare `.Value = "249.25"` assignment would let Excel auto-convert
# that into a real number, which would silently change the cell type. So we
# force the Text number format before assigning, then clear the formatting
# change back off again (ClearFormats) so the cell keeps its original style
# while its content stays textual.
$numericTextUpdates = @{
    "D2" = "249.25"
    "G2" = "14"
    "D3" = "24.01"
    "G3" = "14"
    "D4" = "5.958"
    "G4" = "14"
    "D5" = "0.05903"
    "G5" = "14"
    "D6" = "3.424"
    "G6" = "14"
    "G7" = "14"
    "D8" = "1.332"
    "G8" = "14"
    "D9" = "0.7974"
    "G9" = "14"
    "D10" = "0.1474"
    "G10" = "14"
    "D11" = "0.07830"
    "G11" = "14"
    "D12" = "0.03324"
    "G12" = "14"
    "D13" = "0.03031"
    "G13" = "14"
    "G14" = "14"
    "D15" = "3.566"
    "G15" = "14"
    "D16" = "0.001666"
    "G16" = "14"
    "D17" = "0.04774"
    "G17" = "14"
    "D18" = "0.0006079"
    "G18" = "14"
    "D19" = "0.006231"
    "G19" = "14"
    "D20" = "0.005588"
    "G20" = "14"
    "D21" = "0.001069"
    "G21" = "14"
    "D22" = "0.0001502"
    "G22" = "14"
    "D23" = "3.708"
    "G23" = "14"
    "D24" = "2.212"
    "G24" = "14"
    "D25" = "0.3334"
    "G25" = "14"
    "D26" = "0.1254"
    "G26" = "14"
    "D27" = "0.0006484"
    "G27" = "14"
    "G28" = "14"
    "G29" = "14"
    "G30" = "14"
    "G31" = "14"
    "G32" = "14"
    "G33" = "14"
    "G34" = "14"
    "G35" = "14"
    "G36" = "14"
    "G37" = "14"
    "G38" = "14"
    "G39" = "14"
    "D40" = "0.04409"
    "G40" = "14"
    "D41" = "0.007023"
    "G41" = "14"
    "D42" = "0.003605"
    "G42" = "14"
    "D43" = "0.1066"
    "G43" = "14"
    "D44" = "0.009136"
    "G44" = "14"
    "D45" = "0.002464"
    "G45" = "14"
    "D46" = "0.00005901"
    "G46" = "14"
    "G47" = "14"
    "D48" = "0.9915"
    "G48" = "14"
    "D49" = "0.09911"
    "G49" = "14"
    "D50" = "0.00002103"
    "G50" = "14"
    "G51" = "14"
}
foreach ($ref in $numericTextUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $numericTextUpdates[$ref]
    $cell.ClearFormats()
}

# These two cells hold ordinary (non-numeric-looking) text, so a direct
# assignment is safe and keeps the inline-string type as-is.
$textUpdates = @{
    "E18" = "17OneONE"
    "E51" = "50SpecialPowerGoldSPGWorstin24h"
}
foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}
